$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.370383977890015
$ws.Range("B1").Value = 2.531894445419312
$ws.Range("C1").Value = 2.00196385383606
$ws.Range("D1").Value = 1.968241214752197
$ws.Range("E1").Value = 2.156696081161499
